$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Merge the split "Trocar de usuário: " runs back into a single run.
#    (Find/Replace with identical text causes the matched runs to be
#    merged into one run using the formatting of the first run.)
# ------------------------------------------------------------------
$find = $d.Content.Find
$find.Execute("Trocar de usu" + [char]0xE1 + "rio: ", $false, $false, $false, $false, $false, `
               $true, 1, $false, "Trocar de usu" + [char]0xE1 + "rio: ", 2) | Out-Null

# ------------------------------------------------------------------
# 2) Add <w:noProof/> to every run that hosts an inline picture
#    (Word stamps this onto image runs whenever it resaves a doc).
# ------------------------------------------------------------------
for ($i = 1; $i -le $d.InlineShapes.Count; $i++) {
    $d.InlineShapes.Item($i).Range.NoProofing = 1
}

# ------------------------------------------------------------------
# 3) Turn the two empty paragraphs right after "rm" into:
#       - a centered / underlined heading "Não consegui executar o comando"
#       - a following blank paragraph that keeps an underlined paragraph
#         mark (so the underline "carries" visually to the empty line)
# ------------------------------------------------------------------
$titulo = "N" + [char]0xE3 + "o consegui executar o comando"

$p1 = $d.Paragraphs.Item(37)
$r1 = $p1.Range
$r1.Text = $titulo
$r1b = $p1.Range
$r1b.Font.Name = "Arial"
$r1b.Font.NameAscii = "Arial"
$r1b.Font.NameBi = "Arial"
$r1b.Font.Size = 14
$r1b.Font.SizeBi = 14
$r1b.Font.Underline = 1
$p1.Format.Alignment = 1

$p2 = $d.Paragraphs.Item(38)
$r2 = $p2.Range
$r2.Text = "X"
$r2b = $p2.Range
$r2b.Underline = 1
$r2c = $d.Range($r2b.Start, $r2b.Start + 1)
$r2c.Text = ""

# ------------------------------------------------------------------
# 4) Same transformation for the empty paragraph right after "touch"
#    (no following blank paragraph needs touching this time).
# ------------------------------------------------------------------
$p3 = $d.Paragraphs.Item(40)
$r3 = $p3.Range
$r3.Text = $titulo
$r3b = $p3.Range
$r3b.Font.Name = "Arial"
$r3b.Font.NameAscii = "Arial"
$r3b.Font.NameBi = "Arial"
$r3b.Font.Size = 14
$r3b.Font.SizeBi = 14
$r3b.Font.Underline = 1
$p3.Format.Alignment = 1

Write-Output "done"
